# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#  - Row 2 (fe289fc4-...md) and Row 3 (ffff6b28143b-...md) are now in sync
#    with en-US, so their Status moves from "Ready for handoff" to
#    "Handed back: in sync with en-US".
#  - Their "Latest Target File" (E) / "Latest Handback File" (F) columns
#    get populated with the md / xlf files that were handed back.
#  - Their "Latest Handback DateTime" (G) gets stamped with the real
#    handback time instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Locale-sheet metadata: sheet name, md hyperlink target, xlf hyperlink
# target, and the handback timestamp text for that locale.
$locales = @(
    @{
        Sheet = "zh-cn"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/50e60796851a1e2ef110a6956743292511a9f838/e2e/fe289fc4-9e55-4194-a23f-dd16b44051a8.md"
        MdName = "fe289fc4-9e55-4194-a23f-dd16b44051a8.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c652741f2856434d3dc1c658b9de70c269305b9f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.zh-cn.xlf"
        XlfName = "fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.zh-cn.xlf"
        HandbackDateTime = "2016-02-17 06:35:29"
    },
    @{
        Sheet = "de-de"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/50e60796851a1e2ef110a6956743292511a9f838/e2e/fe289fc4-9e55-4194-a23f-dd16b44051a8.md"
        MdName = "fe289fc4-9e55-4194-a23f-dd16b44051a8.md"
        XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0fb1afe1a0bfc1f1b65579a20f0924799279d92a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.de-de.xlf"
        XlfName = "fe289fc4-9e55-4194-a23f-dd16b44051a8.0cec61a9d3e388cacdf259b410cc67ac11c93aed.de-de.xlf"
        HandbackDateTime = "2016-02-17 06:35:47"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Row 2 and Row 3 are the two files that are now handed back / in sync.
    foreach ($row in 2, 3) {
        $ws.Cells.Item($row, 5).Value = $locale.MdName            # E: Latest Target File
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $locale.MdUrl, "", "", $locale.MdName)

        $ws.Cells.Item($row, 6).Value = $locale.XlfName           # F: Latest Handback File
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $locale.XlfUrl, "", "", $locale.XlfName)

        $ws.Cells.Item($row, 7).Value = $locale.HandbackDateTime  # G: Latest Handback DateTime
    }
}

# Update the Status column text everywhere it currently reads
# "Ready for handoff" (Overview sheet + both locale sheets), so the shared
# string itself flips to the new status for every row that references it.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldStatus, $newStatus)
}
